# Update the 丽水-漫展信息 workbook: refresh the exhibition list for the
# week of 2024-08-10 through 2024-09-16, dropping the two 2024-08-03
# events and shifting everything else up so the sheets end with 4 data
# rows (A1:I5) instead of 6 (A1:I7).
#
# Sheets "展览" (exhibitions) and "全部类型" (all types) carry the same
# table and both need the identical update; "演出" / "本地生活" only
# contain a header row and are left untouched.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- drop the trailing two rows (old rows 6 & 7) so the used range
    #     shrinks from A1:I7 down to A1:I5 ---
    $ws.Range("A6:A7").EntireRow.Delete()

    # Force the date-like "开始时间" column to stay plain text instead of
    # being auto-coerced into a date serial by the Value assignment.
    $ws.Range("B2:B5").NumberFormat = "@"

    # --- row 2: 丽水·CCAC动漫七夕（回馈展） ---
    $ws.Cells.Item(2, 1).Value = 1
    $ws.Cells.Item(2, 2).Value = "2024-08-10"
    $ws.Cells.Item(2, 3).Value = "丽水·CCAC动漫七夕（回馈展）"
    $ws.Cells.Item(2, 4).Value = "中东路848号(解放街交汇) 飞达国际大酒店"
    $ws.Cells.Item(2, 5).Value = "2024.08.10 09:00-08.10 17:00"
    $ws.Cells.Item(2, 6).Value = 100
    $ws.Cells.Item(2, 7).Value = 29.9
    $ws.Cells.Item(2, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86567"
    $ws.Cells.Item(2, 9).Value = "//i0.hdslb.com/bfs/openplatform/202405/tsOzbBRx1717015539538.png"

    # --- row 3: 丽水·AEO纯白礼赞动漫嘉年华 ---
    $ws.Cells.Item(3, 1).Value = 2
    $ws.Cells.Item(3, 2).Value = "2024-08-17"
    $ws.Cells.Item(3, 3).Value = "丽水·AEO纯白礼赞动漫嘉年华"
    $ws.Cells.Item(3, 4).Value = "城北街1001号 爱依·时尚婚宴中心"
    $ws.Cells.Item(3, 5).Value = "2024.08.17 09:00-08.17 18:00"
    $ws.Cells.Item(3, 6).Value = 918
    $ws.Cells.Item(3, 7).Value = 65
    $ws.Cells.Item(3, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86779"
    $ws.Cells.Item(3, 9).Value = "//i2.hdslb.com/bfs/openplatform/202406/MxJ3oNjt1717405405850.jpeg"

    # --- row 4: 丽水·R动漫嘉年华 ---
    $ws.Cells.Item(4, 1).Value = 3
    $ws.Cells.Item(4, 2).Value = "2024-08-24"
    $ws.Cells.Item(4, 3).Value = "丽水·R动漫嘉年华"
    $ws.Cells.Item(4, 4).Value = "中东路848号(解放街交汇) 飞达国际大酒店"
    $ws.Cells.Item(4, 5).Value = "2024.08.24 09:30-08.24 17:00"
    $ws.Cells.Item(4, 6).Value = 216
    $ws.Cells.Item(4, 7).Value = 45
    $ws.Cells.Item(4, 8).Value = "https://show.bilibili.com/platform/detail.html?id=89651"
    $ws.Cells.Item(4, 9).Value = "//i0.hdslb.com/bfs/openplatform/202407/7o5ALbAM1721383424201.jpeg"

    # --- row 5: 丽水·LZ栗子动漫游戏嘉年华 ---
    $ws.Cells.Item(5, 1).Value = 4
    $ws.Cells.Item(5, 2).Value = "2024-09-16"
    $ws.Cells.Item(5, 3).Value = "丽水·LZ栗子动漫游戏嘉年华"
    $ws.Cells.Item(5, 4).Value = "城北街798号 莱茵体育生活馆"
    $ws.Cells.Item(5, 5).Value = "2024.09.16 09:30-09.16 17:00"
    $ws.Cells.Item(5, 6).Value = 427
    $ws.Cells.Item(5, 7).Value = 65
    $ws.Cells.Item(5, 8).Value = "https://show.bilibili.com/platform/detail.html?id=87480"
    $ws.Cells.Item(5, 9).Value = "//i1.hdslb.com/bfs/openplatform/202406/bATqcZhH1719285865931.jpeg"
}
